$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must remain plain text
# (matching the original inline-string cells), so force a temporary Text
# number format while assigning the value, then restore the default style
# so no stray formatting is left behind.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.225.85'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.219.64'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.41%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '298.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '89.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.560'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.01'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.489'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.78'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0780'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.94'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.16%  '
$ws.Range('E13').Value = '  -0.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.559.90'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.212.90'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.777'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '44.009.73'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0908'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.88'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -8.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -10.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '64.47'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.05'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.79'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.35%  '
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.88%  '
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.14'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.41'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '148.56'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.36'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -10.55%  '
$ws.Range('E33').Value = '  -3.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0746'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.115'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.102'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.79'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -12.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.67'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.56%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0301'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.16'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.50'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.60%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.01'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('B43').Value = 'Celestia'
$ws.Range('C43').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.04'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -10.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.810.01'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.76'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.94%  '
$ws.Range('E46').Value = '  -8.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '73.95'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.70%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '93.64'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.53%  '
$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '66.27'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.73%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '13.75'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.442.23'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.30%  '
